# Applies updated crypto price/volume data to Sheet1 (per Aug 29 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.108.42'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.651.21'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5214'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2645'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.89%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06332'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07684'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.60%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.619'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.18%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.661.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.25%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.880.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5593'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.86%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8150'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.71%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.110.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.625'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.82%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '190.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.82%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.931'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.36%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.94%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1189'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.29%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.215'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("E29").Value = '  +2.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05466'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.98%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.269'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.443'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.26%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.361'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.558'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.11%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9494'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.784'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.55%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.401'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.39%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5630'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.47%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01575'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.86%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.857'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.73%  '

# Row 41
$ws.Range("E41").Value = '  -0.21%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8322'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.07%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.028.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.69%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.791.35'
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.51%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈108'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.11%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9988'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.69%  '

# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4339'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.23%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.996'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.51%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05173'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.98%  '
